# Update the crypto price/volume table (rows 2-51) to the latest scrape.
# Row 22 and 23 also swap Coin/Link (Avalanche now ranks above Uniswap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.682.95'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').Value = '1.638.26'
$ws.Range('E3').Value = '  -0.65%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.67'
$ws.Range('D5').Style = "Normal"

$ws.Range('E6').Value = '  -2.26%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '23.08'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.13%  '

$ws.Range('E9').Value = '  +0.71%  '

$ws.Range('E10').Value = '  -0.11%  '

$ws.Range('E11').Value = '  -0.13%  '

$ws.Range('D12').Value = '1.869.84'
$ws.Range('E12').Value = '  -0.66%  '

$ws.Range('D13').Value = '1.632.73'
$ws.Range('E13').Value = '  -0.92%  '

$ws.Range('E14').Value = '  +0.30%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.562'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.37%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.69'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.18%  '

$ws.Range('D17').Value = '27.656.79'
$ws.Range('E17').Value = '  +0.42%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '230.29'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.54%  '

$ws.Range('E19').Value = '  +2.07%  '

$ws.Range('D20').Value = '0.0₃0723'
$ws.Range('E20').Value = '  -0.58%  '

$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.25'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.83%  '

$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.31'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.45%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.06'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.47%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '150.71'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.25%  '

$ws.Range('E26').Value = '  -1.08%  '

$ws.Range('E27').Value = '  -1.12%  '

$ws.Range('E28').Value = '  -0.07%  '

$ws.Range('E29').Value = '  -0.19%  '

$ws.Range('E30').Value = '  +0.27%  '

$ws.Range('E31').Value = '  -0.25%  '

$ws.Range('E32').Value = '  -0.13%  '

$ws.Range('D33').Value = '1.458.54'
$ws.Range('E33').Value = '  +1.72%  '

$ws.Range('E34').Value = '  -2.74%  '

$ws.Range('E35').Value = '  -2.16%  '

$ws.Range('E36').Value = '  -0.15%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.567'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.39%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.880'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.71%  '

$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.893'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +8.74%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '69.19'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.73%  '

$ws.Range('E42').Value = '  -0.04%  '

$ws.Range('E43').Value = '  -1.23%  '

$ws.Range('E44').Value = '  +0.57%  '

$ws.Range('E45').Value = '  -0.51%  '

$ws.Range('E46').Value = '  -0.80%  '

$ws.Range('D47').Value = '1.779.69'
$ws.Range('E47').Value = '  -0.72%  '

$ws.Range('E48').Value = '  +2.58%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '86.68'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.60%  '

$ws.Range('E50').Value = '  +0.20%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.81'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.33%  '
